$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.666448
$ws.Range("H2").Value = 3.332896
$ws.Range("I2").Value = 0.08698774157534103
$ws.Range("J2").Value = 0.0653244391585683
$ws.Range("Q2").Value = 0.05708639817066666
$ws.Range("R2").Value = 0.342518389024
$ws.Range("S2").Value = 0.08698774157534103
$ws.Range("T2").Value = 0.0653244391585683

# Row 3
$ws.Range("I3").Value = 0.1669502665149541
$ws.Range("J3").Value = 0.1880598173367416
$ws.Range("S3").Value = 0.1669502665149541
$ws.Range("T3").Value = 0.1880598173367416

# Row 4
$ws.Range("G4").Value = 2.617047
$ws.Range("H4").Value = 7.851141
$ws.Range("I4").Value = 0.1366085279147753
$ws.Range("J4").Value = 0.1538816040404024
$ws.Range("Q4").Value = 0.089650434381
$ws.Range("R4").Value = 0.806853909429
$ws.Range("S4").Value = 0.1366085279147753
$ws.Range("T4").Value = 0.1538816040404024

# Row 5
$ws.Range("G5").Value = 4.784714
$ws.Range("H5").Value = 9.569428
$ws.Range("I5").Value = 0.2497596474320929
$ws.Range("J5").Value = 0.1875598630045162
$ws.Range("Q5").Value = 0.1639067576886667
$ws.Range("R5").Value = 0.983440546132
$ws.Range("S5").Value = 0.2497596474320929
$ws.Range("T5").Value = 0.1875598630045162

# Row 6
$ws.Range("G6").Value = 4.899255333333334
$ws.Range("H6").Value = 14.697766
$ws.Range("I6").Value = 0.2557386470190557
$ws.Range("J6").Value = 0.2880747916628283
$ws.Range("Q6").Value = 0.1678305237837778
$ws.Range("R6").Value = 1.510474714054
$ws.Range("S6").Value = 0.2557386470190557
$ws.Range("T6").Value = 0.2880747916628283

# Row 7
$ws.Range("G7").Value = 1.991497666666667
$ws.Range("H7").Value = 5.974493
$ws.Range("I7").Value = 0.1039551695437809
$ws.Range("J7").Value = 0.117099484796943
$ws.Range("Q7").Value = 0.06822140790188888
$ws.Range("R7").Value = 0.613992671117
$ws.Range("S7").Value = 0.1039551695437809
$ws.Range("T7").Value = 0.117099484796943
